$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13795
$ws1.Range("F9").Value = 13893
$ws1.Range("F10").Value = 14755
$ws1.Range("F11").Value = 1
$ws1.Range("F26").Value = 5728
$ws1.Range("F30").Value = 48
$ws1.Range("F32").Value = 260

# Sheet "全部类型" (All types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13795
$ws4.Range("F10").Value = 13893
$ws4.Range("F11").Value = 14755
$ws4.Range("F12").Value = 1
$ws4.Range("F27").Value = 5728
$ws4.Range("F31").Value = 48
$ws4.Range("F33").Value = 260
